$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 23.5
$ws.Range("I42").Value = 22
$ws.Range("J42").Value = 25
$ws.Range("K42").Value = 66
$ws.Range("L42").Value = 75
$ws.Range("M42").Value = 164
$ws.Range("N42").Value = -535

$ws.Range("H131").Value = 100000000
$ws.Range("I131").Value = 100000000
$ws.Range("K131").Value = 300000000
$ws.Range("M131").Value = -299994960

$ws.Range("H138").Value = 1760.3143
$ws.Range("I138").Value = 1423.0416
$ws.Range("J138").Value = 2496.182
$ws.Range("K138").Value = 4269.1248
$ws.Range("L138").Value = 7488.545999999999
$ws.Range("M138").Value = 870.8752000000004
$ws.Range("N138").Value = -17768.546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 734.75
$ws.Range("J4").Value = 734.75
$ws.Range("L4").Value = 734.75
$ws.Range("N4").Value = -966.75

$ws.Range("H122").Value = 2235.4
$ws.Range("I122").Value = 1700.5
$ws.Range("K122").Value = 5101.5
$ws.Range("M122").Value = -2651.5

$ws.Range("H132").Value = 2098.6765
$ws.Range("I132").Value = 1691.1
$ws.Range("J132").Value = 5155.5
$ws.Range("K132").Value = 5073.299999999999
$ws.Range("L132").Value = 15466.5
$ws.Range("M132").Value = -2543.299999999999
$ws.Range("N132").Value = -20526.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3675.15
$ws.Range("I86").Value = 3743.9285
$ws.Range("J86").Value = 3514.6667
$ws.Range("K86").Value = 3743.9285
$ws.Range("L86").Value = 3514.6667
$ws.Range("M86").Value = -2620.9285
$ws.Range("N86").Value = -5760.6667

$ws.Range("H89").Value = 3675.15
$ws.Range("I89").Value = 3743.9285
$ws.Range("J89").Value = 3514.6667
$ws.Range("K89").Value = 18719.6425
$ws.Range("L89").Value = 17573.3335
$ws.Range("M89").Value = -13103.6425
$ws.Range("N89").Value = -28805.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("K6").Value = 1
$ws.Range("M6").Value = 112

$ws.Range("H16").Value = 2024.2858
$ws.Range("I16").Value = 1702.6666
$ws.Range("J16").Value = 2603.2
$ws.Range("K16").Value = 1702.6666
$ws.Range("L16").Value = 2603.2
$ws.Range("M16").Value = -1415.6666
$ws.Range("N16").Value = -3177.2

$ws.Range("H22").Value = 458.05884
$ws.Range("I22").Value = 288.42856
$ws.Range("J22").Value = 1249.6666
$ws.Range("K22").Value = 288.42856
$ws.Range("L22").Value = 1249.6666
$ws.Range("M22").Value = 61.57144
$ws.Range("N22").Value = -1949.6666

$ws.Range("H107").Value = 1582.0416
$ws.Range("J107").Value = 1949.5
$ws.Range("L107").Value = 1949.5
$ws.Range("N107").Value = -5789.5

$ws.Range("H113").Value = 2024.2858
$ws.Range("I113").Value = 1702.6666
$ws.Range("J113").Value = 2603.2
$ws.Range("K113").Value = 1702.6666
$ws.Range("L113").Value = 2603.2
$ws.Range("M113").Value = 467.3334
$ws.Range("N113").Value = -6943.2

$ws.Range("H132").Value = 2676064.5
$ws.Range("I132").Value = 3032456.5
$ws.Range("K132").Value = 9097369.5
$ws.Range("M132").Value = -9094839.5

$ws.Range("H134").Value = 2788376.5
$ws.Range("I134").Value = 3574127.2
$ws.Range("J134").Value = 169207.67
$ws.Range("K134").Value = 10722381.6
$ws.Range("L134").Value = 507623.01
$ws.Range("M134").Value = -10719846.6
$ws.Range("N134").Value = -512693.01

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1021.4545
$ws.Range("I5").Value = 956
$ws.Range("K5").Value = 2868
$ws.Range("M5").Value = -2756

$ws.Range("H23").Value = 113691.336
$ws.Range("J23").Value = 113691.336
$ws.Range("L23").Value = 341074.008
$ws.Range("N23").Value = -341544.008

$ws.Range("H40").Value = 73.75
$ws.Range("J40").Value = 100
$ws.Range("L40").Value = 400
$ws.Range("N40").Value = -538

$ws.Range("H68").Value = 104409.3
$ws.Range("I68").Value = 4166.6665
$ws.Range("J68").Value = 147370.42
$ws.Range("K68").Value = 12499.9995
$ws.Range("L68").Value = 442111.26
$ws.Range("M68").Value = -11688.9995
$ws.Range("N68").Value = -443733.26

$ws.Range("H71").Value = 104409.3
$ws.Range("I71").Value = 4166.6665
$ws.Range("J71").Value = 147370.42
$ws.Range("K71").Value = 37499.9985
$ws.Range("L71").Value = 1326333.78
$ws.Range("M71").Value = -33443.9985
$ws.Range("N71").Value = -1334445.78

$ws.Range("H135").Value = 1021.4545
$ws.Range("I135").Value = 956
$ws.Range("K135").Value = 8604
$ws.Range("M135").Value = -6069

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 2296
$ws.Range("J3").Value = 2296
$ws.Range("L3").Value = 2296
$ws.Range("N3").Value = -2528

$ws.Range("H10").Value = 2503899.5
$ws.Range("J10").Value = 7799
$ws.Range("L10").Value = 7799
$ws.Range("N10").Value = -8137

$ws.Range("H132").Value = 4290.5557
$ws.Range("I132").Value = 3580
$ws.Range("K132").Value = 10740
$ws.Range("M132").Value = -8210

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1821.5
$ws.Range("I61").Value = 643
$ws.Range("K61").Value = 643
$ws.Range("M61").Value = -441

$ws.Range("H100").Value = 10502.125
$ws.Range("I100").Value = 10742.267
$ws.Range("K100").Value = 10742.267
$ws.Range("M100").Value = -10201.267

$ws.Range("H113").Value = 1821.5
$ws.Range("I113").Value = 643
$ws.Range("K113").Value = 643
$ws.Range("M113").Value = 1527

$ws.Range("H132").Value = 2686.7917
$ws.Range("I132").Value = 1995.2858
$ws.Range("J132").Value = 3654.9
$ws.Range("K132").Value = 5985.857400000001
$ws.Range("L132").Value = 10964.7
$ws.Range("M132").Value = -3455.857400000001
$ws.Range("N132").Value = -16024.7

$ws.Range("H136").Value = 2122.4666
$ws.Range("I136").Value = 1778.6666
$ws.Range("J136").Value = 3497.6667
$ws.Range("K136").Value = 5335.9998
$ws.Range("L136").Value = 10493.0001
$ws.Range("M136").Value = -2785.9998
$ws.Range("N136").Value = -15593.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 53999
$ws.Range("I54").Value = 53999
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 53999
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -53479
$ws.Range("N54").ClearContents()

$ws.Range("H81").Value = 2004.7
$ws.Range("I81").Value = 1068.375
$ws.Range("J81").Value = 5750
$ws.Range("K81").Value = 2136.75
$ws.Range("L81").Value = 11500
$ws.Range("M81").Value = -1075.75
$ws.Range("N81").Value = -13622

$ws.Range("H84").Value = 2004.7
$ws.Range("I84").Value = 1068.375
$ws.Range("J84").Value = 5750
$ws.Range("K84").Value = 10683.75
$ws.Range("L84").Value = 57500
$ws.Range("M84").Value = -5379.75
$ws.Range("N84").Value = -68108

$ws.Range("H132").Value = 2746.5789
$ws.Range("I132").Value = 2598.9167
$ws.Range("J132").Value = 2999.7144
$ws.Range("K132").Value = 7796.750100000001
$ws.Range("L132").Value = 8999.143199999999
$ws.Range("M132").Value = -5266.750100000001
$ws.Range("N132").Value = -14059.1432

Write-Output "Applied scheduled profit-sheet updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR"
